# Daily attendance processing - 2025-10-29 08:53:33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column I (9th column) width: 14 -> 10
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 9.17

# ---------------------------------------------------------------------------
# 2. "Recorded By" (column G) values: swap the first two comma-separated
#    entries (re-ordering the recorder names) for the affected rows.
# ---------------------------------------------------------------------------
$gSwapRows = @(2,3,4,6,7,10,12,13,14,15,18,19,20,21,22,24,29,30,31,33,34,37,39,40,41,42,45,46,47,48,49,51,56,57,58,60,61,64,66,67,68,69,72,73,74,75,76,78,86,87,88,89,90,93,95,102,112,113,114,115,116,119,121,128,138,139,140,141,142,145,147,154)

foreach ($r in $gSwapRows) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    $parts = $val -split ", "
    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp
    $cell.Value2 = ($parts -join ", ")
}

# ---------------------------------------------------------------------------
# Helper: write a plain TEXT value into a cell without Excel's "smart" input
# parser re-interpreting a "NN.N%"-shaped string as a numeric percentage
# (which would both change the stored type and pull in a new numFmt'd
# style). Routing the text through a formula result + paste-values keeps
# the original General-format style (e.g. "s=4") untouched.
# ---------------------------------------------------------------------------
function Set-TextValue($row, $col, $text) {
    $scratch = $ws.Cells.Item(300, 1)
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $target = $ws.Cells.Item($row, $col)
    $target.PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# 3. Recalculated summary-statistics cells (Class Statistics block)
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 12).Value2 = 138
$ws.Cells.Item(7, 12).Value2 = 0
Set-TextValue 9 12 "86.8%"
Set-TextValue 10 12 "71.1%"

# ---------------------------------------------------------------------------
# 4. Per-group rollup rows (18, 19, 20) - Missed/Excused counts & percentages
# ---------------------------------------------------------------------------
$ws.Cells.Item(18, 15).Value2 = 22
$ws.Cells.Item(18, 16).Value2 = 0
Set-TextValue 18 18 "84.6%"
Set-TextValue 18 19 "75.9%"

$ws.Cells.Item(19, 15).Value2 = 22
$ws.Cells.Item(19, 16).Value2 = 0
Set-TextValue 19 18 "84.6%"
Set-TextValue 19 19 "73.6%"

$ws.Cells.Item(20, 15).Value2 = 22
$ws.Cells.Item(20, 16).Value2 = 0
Set-TextValue 20 18 "84.6%"
Set-TextValue 20 19 "81.0%"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Newly-recorded sessions: rows 104 (B2D), 130 (B2E), 156 (B2F)
#    These flip from the "Not Recorded" (pink) style to the "Recorded"
#    (green) style used by the normal data rows, and gain attendance data.
# ---------------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A104:I104").PasteSpecial(-4122)
$ws.Range("A130:I130").PasteSpecial(-4122)
$ws.Range("A156:I156").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(104, 7).Value2 = "dnasr281@gmail.com"
$ws.Cells.Item(104, 8).Value2 = "35/56"
$ws.Cells.Item(104, 9).Value2 = "Recorded"

$ws.Cells.Item(130, 7).Value2 = "dnasr281@gmail.com"
$ws.Cells.Item(130, 8).Value2 = "33/55"
$ws.Cells.Item(130, 9).Value2 = "Recorded"

$ws.Cells.Item(156, 7).Value2 = "dnasr281@gmail.com"
$ws.Cells.Item(156, 8).Value2 = "40/57"
$ws.Cells.Item(156, 9).Value2 = "Recorded"
